# Extend the "17.4.1 Debt service" table with two more years (2023, 2024):
# columns T (20) and U (21) get appended, mirroring the formatting of the
# preceding column S (19) for rows 3, 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 3 (thin divider row under the headers): extend the bottom border
# formatting two more columns (empty cells, same style as S3) ---
$ws.Range("S3").Copy()
$ws.Range("T3:U3").PasteSpecial(-4122)

# --- row 4 (year header row): extend with 2023 / 2024 ---
$ws.Range("S4").Copy()
$ws.Range("T4:U4").PasteSpecial(-4122)
$ws.Cells.Item(4, 20).Value = 2023
$ws.Cells.Item(4, 21).Value = 2024

# --- row 5 (data row): extend with the new data points ---
$ws.Range("S5").Copy()
$ws.Range("T5:U5").PasteSpecial(-4122)
$ws.Cells.Item(5, 20).Value = 10.8
$ws.Cells.Item(5, 21).Value = 6.5

$ws.Application.CutCopyMode = 0

# row 5 grew taller to fit the (wrapped) text now that the table is wider
$ws.Rows.Item(5).RowHeight = 41.25

# the newly-used columns D:U get an explicit (default) width, matching the
# rest of the numeric columns in the sheet
$ws.Range("D1:U1").ColumnWidth = 7.9

# reset the selection back to the sheet's home cell
[void]$ws.Range("A1").Select()
